$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting existing rows 9-107 down to 10-108
$ws.Rows.Item(9).Insert()

# Populate the new row 9 with data
$ws.Cells.Item(9, 1).Value = 3
$ws.Cells.Item(9, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(9, 3).Value = "Coquimbo"
$ws.Cells.Item(9, 4).Value = 44537
$ws.Cells.Item(9, 5).Value = 5
$ws.Cells.Item(9, 6).Value = 100112030
$ws.Cells.Item(9, 7).Value = "Poroto granado"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 80
$ws.Cells.Item(9, 11).Value = 45000
$ws.Cells.Item(9, 12).Value = 46000
$ws.Cells.Item(9, 13).Value = 45500
$ws.Cells.Item(9, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(9, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(9, 16).Value = 1820
$ws.Cells.Item(9, 17).Value = 25
$ws.Cells.Item(9, 18).Value = "Hortaliza"
